$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1) -----------------------------------------------
$ws.Cells.Item(1, 2).Value  = "species"
$ws.Cells.Item(1, 3).Value  = "debtor"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "total"
$ws.Cells.Item(1, 6).Value  = "register_date"
$ws.Cells.Item(1, 7).Value  = "register_reason"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 (104 - 土地房屋貸款 / 合作金庫商業銀行) ----------------------
$ws.Cells.Item(2, 2).Value  = "土地房屋貸款"
$ws.Cells.Item(2, 3).Value  = "徐耀昌"
$ws.Cells.Item(2, 4).Value  = "合作金庫商業銀行"
$ws.Cells.Item(2, 5).Value  = 15279280
$ws.Cells.Item(2, 6).Value  = "99年01月19曰"
$ws.Cells.Item(2, 7).Value  = "抵押貸款"
$ws.Cells.Item(2, 8).Value  = "debt"
$ws.Cells.Item(2, 9).Value  = "normal"
$ws.Cells.Item(2, 10).Value = "2012-04-24"
$ws.Cells.Item(2, 11).Value = "徐耀昌"
$ws.Cells.Item(2, 12).Value = 921
$ws.Cells.Item(2, 13).Value = "tmp6e501"
$ws.Cells.Item(2, 14).Value = 104

# --- Row 3 (105 - 房M貸款 / 臺灣土地銀行) -------------------------------
$ws.Cells.Item(3, 2).Value  = "房M貸款"
$ws.Cells.Item(3, 3).Value  = "徐耀昌"
$ws.Cells.Item(3, 4).Value  = "臺灣土地銀行"
$ws.Cells.Item(3, 5).Value  = 12000000
$ws.Cells.Item(3, 6).Value  = "100年02月11U"
$ws.Cells.Item(3, 7).Value  = "抵押貸款"
$ws.Cells.Item(3, 8).Value  = "debt"
$ws.Cells.Item(3, 9).Value  = "normal"
$ws.Cells.Item(3, 10).Value = "2012-04-24"
$ws.Cells.Item(3, 11).Value = "徐耀昌"
$ws.Cells.Item(3, 12).Value = 921
$ws.Cells.Item(3, 13).Value = "tmp6e501"
$ws.Cells.Item(3, 14).Value = 105
